$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C ("Scale" column). Use xlFormatFromRightOrBelow so the
# shifted-right old-C formatting lines up; we then repaint column C's format from the
# (now shifted) column D, which carries the original column-C formatting.
$ws.Columns("C").Insert(-4161)

$ws.Range("D1:D24").Copy() | Out-Null
$ws.Range("C1:C24").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Fill in the scale (spatial extent) for each covariate group
$ws.Range("C4").Value = "750 m"
$ws.Range("C5").Value = "750 m"
$ws.Range("C6").Value = "750 m"
$ws.Range("C7").Value = "750 m"
$ws.Range("C8").Value = "750 m"

$ws.Range("C10").Value = "3000 m"
$ws.Range("C11").Value = "3000 m"
$ws.Range("C12").Value = "3000 m"
$ws.Range("C13").Value = "3000 m"
$ws.Range("C14").Value = "3000 m"
$ws.Range("C15").Value = "3000 m"

$ws.Range("C17").Value = "4000 m"
$ws.Range("C18").Value = "4000 m"
$ws.Range("C19").Value = "4000 m"
$ws.Range("C20").Value = "4000 m"
$ws.Range("C21").Value = "4000 m"
$ws.Range("C22").Value = "4000 m"

# Header for the new column (typed last)
$ws.Range("C2").Value = "Scale"

# Normalize the COMP sub-model row labels (B10:B14, B15) onto their existing equivalent
# style (matching B6's / B8's format) so no duplicate style entries remain, same as Excel does
# when the formatting is re-applied / cleaned up on save.
$ws.Range("B6").Copy() | Out-Null
$ws.Range("B10:B14").PasteSpecial(-4122) | Out-Null
$ws.Range("B8").Copy() | Out-Null
$ws.Range("B15").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Column widths: widen B slightly, narrow the new Scale column C
$ws.Columns("B").ColumnWidth = 23.417
$ws.Columns("C").ColumnWidth = 8.417

# View state: zoomed in, selection moved to the new rightmost (Mean) column
$ws.Select()
$excel.ActiveWindow.Zoom = 175
$ws.Range("H30").Select()
